$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34/35: Hedera and InternetComputer(DFINITY) swap places, with updated price/volume data.
$ws.Cells.Item(34, 2).Value = "Hedera"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.0603"
$ws.Cells.Item(34, 5).Value = "  -0.58%  "

$ws.Cells.Item(35, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "4.50"
$ws.Cells.Item(35, 5).Value = "  -1.40%  "

# Remaining price (D) and volume (E) updates for all other rows.
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "37.779.14"
$ws.Cells.Item(2, 5).Value = "  -1.09%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.037.78"
$ws.Cells.Item(3, 5).Value = "  -0.92%  "
$ws.Cells.Item(4, 5).Value = "  +0.12%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "227.96"
$ws.Cells.Item(5, 5).Value = "  -0.16%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.607"
$ws.Cells.Item(6, 5).Value = "  -1.54%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "60.55"
$ws.Cells.Item(7, 5).Value = "  -0.63%  "
$ws.Cells.Item(8, 5).Value = "  +0.12%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.375"
$ws.Cells.Item(9, 5).Value = "  -2.59%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.0825"
$ws.Cells.Item(10, 5).Value = "  -0.03%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.104"
$ws.Cells.Item(11, 5).Value = "  +0.24%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "2.347.64"
$ws.Cells.Item(12, 5).Value = "  -0.48%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "14.56"
$ws.Cells.Item(13, 5).Value = "  -1.81%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "21.00"
$ws.Cells.Item(14, 5).Value = "  -0.63%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.775"
$ws.Cells.Item(15, 5).Value = "  +1.67%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "5.22"
$ws.Cells.Item(16, 5).Value = "  -1.63%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "2.049.80"
$ws.Cells.Item(17, 5).Value = "  -0.91%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "37.761.29"
$ws.Cells.Item(18, 5).Value = "  -0.95%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "69.44"
$ws.Cells.Item(19, 5).Value = "  -0.54%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "5.92"
$ws.Cells.Item(20, 5).Value = "  -3.97%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "0.0₃0822"
$ws.Cells.Item(21, 5).Value = "  -1.26%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "223.94"
$ws.Cells.Item(22, 5).Value = "  -0.59%  "
$ws.Cells.Item(23, 5).Value = "  -0.06%  "
$ws.Cells.Item(24, 5).Value = "  -0.39%  "
$ws.Cells.Item(25, 5).Value = "  +2.53%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "168.13"
$ws.Cells.Item(26, 5).Value = "  +0.86%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "9.31"
$ws.Cells.Item(27, 5).Value = "  +0.94%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.129"
$ws.Cells.Item(28, 5).Value = "  -2.03%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "18.77"
$ws.Cells.Item(29, 5).Value = "  -1.20%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.26"
$ws.Cells.Item(30, 5).Value = "  -2.56%  "
$ws.Cells.Item(31, 5).Value = "  -0.73%  "
$ws.Cells.Item(32, 5).Value = "  +7.46%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "4.39"
$ws.Cells.Item(33, 5).Value = "  -2.35%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "6.55"
$ws.Cells.Item(36, 5).Value = "  +3.75%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.34"
$ws.Cells.Item(37, 5).Value = "  +2.10%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "3.43"
$ws.Cells.Item(38, 5).Value = "  +4.82%  "
$ws.Cells.Item(39, 5).Value = "  +0.03%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "18.07"
$ws.Cells.Item(40, 5).Value = "  +7.50%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.533.08"
$ws.Cells.Item(41, 5).Value = "  +0.02%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.0215"
$ws.Cells.Item(42, 5).Value = "  -1.68%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "96.57"
$ws.Cells.Item(43, 5).Value = "  -1.81%  "
$ws.Cells.Item(44, 5).Value = "  -0.67%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.0910"
$ws.Cells.Item(45, 5).Value = "  -2.23%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "4.11"
$ws.Cells.Item(46, 5).Value = "  +2.91%  "
$ws.Cells.Item(47, 5).Value = "  -1.11%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.00"
$ws.Cells.Item(48, 5).Value = "  -1.03%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.94"
$ws.Cells.Item(49, 5).Value = "  -1.63%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "7.06"
$ws.Cells.Item(50, 5).Value = "  +0.03%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "2.233.59"
$ws.Cells.Item(51, 5).Value = "  -0.66%  "
